$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '91.474.06'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.107.62'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.40'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '620.48'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.49%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +8.99%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.369'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.33%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.110.33'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +6.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.202'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.62'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.521.14'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.49'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.698.40'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.154.86'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.70'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.65'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000213'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.84'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '446.68'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.24'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.90'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '91.28'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +8.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.98'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.266.22'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +17.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.246'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +27.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.29'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.171'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +15.85%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.113'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +34.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.82'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +9.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.69'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.19'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +28.42%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '497.31'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.36%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.64'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.30'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.15'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.91'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.696'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '153.84'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.55'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.35'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.35'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.38%  '
